# Spring_High_Data.xlsx edit script
# - Shifts the HH Tide Values in column B down by one row (rows 3-25),
#   pulling a "General"-formatted value into row 3 and giving row 4 the
#   "0.000" number format that row 5 used to have (matches the diff's
#   style churn on B3/B4).
# - Adds a new data row 26 (date 20250227, value 6.512) that used to be
#   an empty placeholder row.
# - Updates the Median (B28) / Average (B29) summary cells to their new
#   values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B3 loses its "0.000" number format (becomes plain/General), B4 gains it ---
$ws.Range("B3").ClearFormats()
$ws.Range("B3").Value = 5

$ws.Range("B4").NumberFormat = "0.000"
$ws.Range("B4").Value = 7.044

# --- Remaining rows: each cell takes on the value that used to be one row below it ---
$ws.Range("B5").Value = 5.82
$ws.Range("B6").Value = 6.434
$ws.Range("B7").Value = 5.981
$ws.Range("B8").Value = 6.496
$ws.Range("B9").Value = 4.078

$ws.Range("B13").Value = 6.191
$ws.Range("B14").Value = 5.974
$ws.Range("B15").Value = 5.951
$ws.Range("B16").Value = 6.155
$ws.Range("B17").Value = 6.49
$ws.Range("B18").Value = 5.256
$ws.Range("B19").Value = 7.543
$ws.Range("B20").Value = 4.734
$ws.Range("B21").Value = 5.745
$ws.Range("B22").Value = 6.217
$ws.Range("B23").Value = 6.155
$ws.Range("B24").Value = 5.581
$ws.Range("B25").Value = 5.945

# --- New row 26: 20250227 / 6.512 (previously a blank filler row) ---
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "20250227"
$ws.Range("A26").ClearFormats()

$ws.Range("B26").NumberFormat = "0.000"
$ws.Range("B26").Value = 6.512

# --- Updated summary stats ---
$ws.Range("B28").Value = 6.009
$ws.Range("B29").Value = 6.044880000000001
